$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data correction: realign row labels with their numeric data, and
#     rename the "unnamed: 1_level_1" header to "total" (commit: correcao
#     nos dados e inicio da analise PNAD 2009) ---

$ws.Range("B2").Value = "total"

$rowData = @(
    @{ Row = 4; Label = "               brasil"; Vals = @(0.6928770063293205,0.7276622813388498,0.7237218662075712,0.9649814785319062,1.069489298763232,2.096332857131479,1.149406663822678) },
    @{ Row = 5; Label = "urbana"; Vals = @(0.8587980411850094,0.891576150167402,0.8970706186869837,1.087885572715446,1.236528651394527,2.107266626755658,1.322686315619079) },
    @{ Row = 6; Label = "rural"; Vals = @(3.25391337464394,3.208105724805562,3.367786289168296,4.170488497400846,3.723584414504003,6.460993726707648,3.823049076506237) },
    @{ Row = 7; Label = "norte"; Vals = @(3.476844813349684,3.539967467924956,3.514237853458653,3.492650478265156,3.799080872373382,6.272179608686657,3.833621214152775) },
    @{ Row = 8; Label = "rondônia"; Vals = @(3.315533400930749,4.066222730551107,3.460686205685528,4.734966341435021,3.673034632728099,8.806889398846353,4.072842100129235) },
    @{ Row = 9; Label = "acre"; Vals = @(10.76603758468494,11.06901070500634,10.67073104353443,14.49689659396446,10.42994617305908,20.62865367213075,10.71026944602688) },
    @{ Row = 10; Label = "amazonas"; Vals = @(4.418281608947436,4.774823161702064,4.359980645584146,7.251866512735738,5.000876783577505,25.73719603333494,4.806420490176979) },
    @{ Row = 11; Label = "roraima"; Vals = @(4.756212922640313,5.325509745814172,5.183421899544729,8.961987330019374,5.485242910712016,17.59729101316955,5.796712616983622) },
    @{ Row = 12; Label = "pará"; Vals = @(8.708223025440608,8.648892917876307,8.935680619516937,8.330555618831855,9.022627375983012,12.06557741173178,9.064163348455002) },
    @{ Row = 13; Label = "amapá"; Vals = @(8.663353080183908,9.61995803518318,8.248005262957353,17.9054132086401,6.762009980762627,24.21046280373702,7.493537321272747) },
    @{ Row = 14; Label = "tocantins"; Vals = @(3.401724653742177,3.287426112114832,3.789953367655787,6.308638321606763,4.898033723522428,13.83587010614151,5.065867695811757) },
    @{ Row = 15; Label = "nordeste"; Vals = @(1.653280897927568,1.70389735416727,1.675108115882935,2.655816658165791,2.036605737637894,4.449161186226763,2.1414008341255) },
    @{ Row = 16; Label = "maranhão"; Vals = @(7.394463318254337,7.684601121422091,7.267828567267048,12.95766891997372,8.758770102825974,16.7482970504253,10.22038145894486) },
    @{ Row = 17; Label = "piauí"; Vals = @(5.896780275014431,6.007528397394296,6.081791979912468,9.453224471453266,7.51465436712801,31.8228074996303,6.694548497980144) },
    @{ Row = 18; Label = "ceará"; Vals = @(4.10944186047277,4.241612002328409,4.191060031579496,4.716903870141481,4.847335349319752,14.0978374864585,4.795045790329499) },
    @{ Row = 19; Label = "rio grande do norte"; Vals = @(9.232800347021263,8.711768661217958,9.870092014996668,11.56204234420452,8.305640431077553,12.33987016934834,8.75531532390947) },
    @{ Row = 20; Label = "paraíba"; Vals = @(5.58828201195286,5.888240376756257,5.544184819289105,6.896499373766177,6.865564114228651,12.04836259437774,6.909562135793559) },
    @{ Row = 21; Label = "pernambuco"; Vals = @(2.843305541804166,3.278372708452765,2.647514661367444,3.523322279187858,4.045369352373674,10.42082215639311,4.059597145290524) },
    @{ Row = 22; Label = "alagoas"; Vals = @(8.573679076448991,8.369221577078443,8.896911676739673,17.72063987643193,9.473575414882262,20.97202805163978,10.38774200215877) },
    @{ Row = 23; Label = "sergipe"; Vals = @(5.924223609358608,6.539279900086489,5.849817669668029,8.053623842991975,7.007288516007485,18.72608251823245,6.97676735608239) },
    @{ Row = 24; Label = "bahia"; Vals = @(2.740742423175238,2.784744104600273,2.870662009988054,5.963876229111223,3.610283934931857,5.899886452971694,3.714714066224043) },
    @{ Row = 25; Label = "sudeste"; Vals = @(0.953118753008872,1.009093433702807,1.042229673345515,1.472035372771579,1.49210702318353,3.029078996907012,1.643598971860037) },
    @{ Row = 26; Label = "minas gerais"; Vals = @(1.894990983554854,2.003310264746617,2.013068495149036,3.20383856773701,2.167821442180994,5.60067134535625,2.250345324204977) },
    @{ Row = 27; Label = "espírito santo"; Vals = @(4.038166149009035,3.851891977863552,4.983357256230466,8.80418932663505,5.227630751158408,12.21895579229582,5.724061329522098) },
    @{ Row = 28; Label = "rio de janeiro"; Vals = @(2.335418274795894,2.382724801554635,2.514937583104198,3.732440796891964,3.156952793980208,5.640401303600758,3.397454989417274) },
    @{ Row = 29; Label = "são paulo"; Vals = @(1.304069382914054,1.406155813004492,1.447354299460135,1.886195956057096,2.738415945606821,4.956486314829799,3.105729685373373) },
    @{ Row = 30; Label = "sul"; Vals = @(1.385022388995266,1.505509081145612,1.446850189662667,1.754023811499629,3.221768414681316,5.087925844676839,3.601614726687955) },
    @{ Row = 31; Label = "paraná"; Vals = @(2.17783554277198,2.327127538423158,2.462165665182943,2.656852338085135,4.90626690032197,8.933270579978922,5.306863704715269) },
    @{ Row = 32; Label = "santa catarina"; Vals = @(3.733299083458498,3.693364198381585,3.961758406306001,4.878412984967039,7.816015427384113,13.16872124476143,9.06162130458714) },
    @{ Row = 33; Label = "rio grande do sul"; Vals = @(1.693546907802575,2.143334308868634,1.541307339217039,1.838887141112,4.9997070252886,7.009568795246629,5.623733054616518) },
    @{ Row = 34; Label = "centro-oeste"; Vals = @(1.405888871601273,1.625597007641775,1.574432473849135,2.594991830296206,2.349123536634424,4.917023236547275,2.451165214796853) },
    @{ Row = 35; Label = "mato grosso do sul"; Vals = @(3.061358952849435,3.688496814466627,3.113295287683079,4.15202663840212,5.298865866465358,10.73614003314723,5.842584303283794) },
    @{ Row = 36; Label = "mato grosso"; Vals = @(2.862124415579702,3.880823947995638,2.593966100580884,8.26198330544913,5.995072884640009,9.593045431716474,6.307604532982513) },
    @{ Row = 37; Label = "goiás"; Vals = @(2.190081720279117,2.324882247299495,2.747284271774886,3.785698292318595,3.320511453848968,8.231513804545303,3.443123915437826) },
    @{ Row = 38; Label = "distrito federal"; Vals = @(3.767948480175316,3.903062668213608,4.000778682071914,5.024768398545165,4.343937480142502,9.378804756806005,4.451807146905777) }
)

foreach ($item in $rowData) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Label
    $vals = $item.Vals
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($r, $c + 2).Value = $vals[$c]
    }
}

# The last two old rows (39 and 40) are no longer needed: their data has
# been absorbed into the realigned rows above, so drop the leftover rows.
$ws.Rows.Item(39).Delete()
$ws.Rows.Item(39).Delete()
